# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (want-to-go count) figures and one event cover image
# across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 508
$ws1.Range("F3").Value  = 6037
$ws1.Range("F5").Value  = 84
$ws1.Range("F6").Value  = 115
$ws1.Range("F8").Value  = 63
$ws1.Range("I8").Value  = "//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg"
$ws1.Range("F9").Value  = 555
$ws1.Range("F10").Value = 36

# --- Sheet "全部类型" --------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 508
$ws4.Range("F3").Value  = 6037
$ws4.Range("F6").Value  = 84
$ws4.Range("F7").Value  = 115
$ws4.Range("F10").Value = 63
$ws4.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202409/hBNwkgri1727595368161.jpeg"
$ws4.Range("F11").Value = 555
$ws4.Range("F12").Value = 36
